$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New quarterly row ("01-07-2021") appended at the bottom of the series.
# Force the date-looking label to be stored as text (matches the other
# "Serie" cells, which are plain shared-string text, not date values),
# then drop the temporary number-format override so the cell keeps the
# workbook's default style like its neighbours above it.
$ws.Range("A59").NumberFormatLocal = "@"
$ws.Range("A59").Value = "01-07-2021"
$ws.Range("A59").ClearFormats()

$ws.Range("B59").Value = 7353
$ws.Range("C59").Value = 497
$ws.Range("D59").Value = 1345
$ws.Range("E59").Value = 993
$ws.Range("F59").Value = 749
$ws.Range("G59").Value = 1370
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2400
$ws.Range("K59").Value = 0
